$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rewrite the "Horarios Oficiales" fragment embedded in the structuredData
#    (column M) for the two listings that remain, switching the old
#    "jueves, De .. a .., Copiar el horario,..." scrape format to the cleaned
#    "Viernes de .. a ..,Sábado de .. a ..,...,Jueves de .. a .." format.
$m2Text = @'

                    <h2>Parque Ecoturístico Vivero Santa Fe</h2>
                    <p><b>Dirección del Parque Ecoturístico: </b>Blvd. Cnel. Enrique Carrola Antuna 909, Ciénega, 34090 Durango, Dgo.</p>
                    <p><b>Teléfono del Parque Ecoturístico: </b>618 235 9375</p>
                    <p><b>Horarios Oficiales: </b>Viernes de 09:30 a 19:30,Sábado de 09:30 a 19:30,Domingo de 09:30 a 19:30,Lunes de 09:30 a 19:30,Martes de 09:30 a 19:30,Miércoles de 09:30 a 19:30,Jueves de 09:30 a 19:30</p>
                    <p><b>Sitio Web: </b>web no disponible</p>
                    <p><b>Ubicación: </b><a href='https://www.google.com.mx/maps/place/Vivero+Santa+Fe/data=!4m7!3m6!1s0x869bb7e17f5ebdfd:0x329bbfbf57717ca0!8m2!3d24.0199748!4d-104.6558555!16s%2Fg%2F11btmr25_0!19sChIJ_b1ef-G3m4YRoHxxV7-_mzI?authuser=0&hl=es&rclk=1'>Mapa del Parque Ecoturístico Vivero Santa Fe</a></p>                        
'@
$ws.Range("M2").Value = $m2Text

$m3Text = @'

                    <h2>Parque Ecoturístico Vivero Las Magnolias</h2>
                    <p><b>Dirección del Parque Ecoturístico: </b>Sauca 103, Jardines de Durango, 34200 Durango, Dgo.</p>
                    <p><b>Teléfono del Parque Ecoturístico: </b>618 129 7673</p>
                    <p><b>Horarios Oficiales: </b>Viernes de 09:00 a 18:00,Sábado de 09:00 a 18:00,Domingo de 09:00 a 17:00,Lunes de 09:00 a 18:00,Martes de 09:00 a 18:00,Miércoles de 09:00 a 18:00,Jueves de 09:00 a 18:00</p>
                    <p><b>Sitio Web: </b>web no disponible</p>
                    <p><b>Ubicación: </b><a href='https://www.google.com.mx/maps/place/Vivero+Las+Magnolias/data=!4m7!3m6!1s0x869bb7b871872b85:0x65cc8cefdc6b7345!8m2!3d24.0374832!4d-104.6351623!16s%2Fg%2F11b6j5b_lw!19sChIJhSuHcbi3m4YRRXNr3O-MzGU?authuser=0&hl=es&rclk=1'>Mapa del Parque Ecoturístico Vivero Las Magnolias</a></p>                        
'@
$ws.Range("M3").Value = $m3Text

# Re-run autofit so the row height stays on "auto" (matches the untouched
# rows) instead of picking up an explicit height from the multi-line value.
$ws.Rows("2:3").EntireRow.AutoFit()

# 2) Drop the last three listings (VIVERO'S AVE DE PARAISO, Viveros Del
#    Guadiana, the duplicate Vivero Santa Fe) - only the header + first two
#    records are kept, shrinking the sheet from A1:O6 down to A1:O3.
$ws.Rows("4:6").Delete()
